# Penalty/Reward system tweak: shift the forecast window forward one week
# and update the corresponding MyForecast values + the dependent Summary
# sheet statistics.
#
# NOTE: several of the literal values being written look like dates or
# plain numbers ("2025-01-12", "47", ...) but the source workbook stores
# them as literal text (t="inlineStr") rather than as real numbers/dates.
# Excel's normal cell-entry parser would otherwise auto-convert these into
# date serials / numeric cells, so a leading apostrophe (quote-prefix) is
# used to force them to stay text, exactly as a user typing '2025-01-12
# into a text cell would.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": Week_Start_Date (B) + MyForecast (D) ---

$ws1.Range("B2").Value  = "'2025-01-12"
$ws1.Range("D2").Value  = 47

$ws1.Range("B3").Value  = "'2025-01-19"
$ws1.Range("D3").Value  = 47

$ws1.Range("B4").Value  = "'2025-01-26"
$ws1.Range("D4").Value  = 48

$ws1.Range("B5").Value  = "'2025-02-02"
$ws1.Range("D5").Value  = 49

$ws1.Range("B6").Value  = "'2025-02-09"
$ws1.Range("D6").Value  = 49

$ws1.Range("B7").Value  = "'2025-02-16"
$ws1.Range("D7").Value  = 50

$ws1.Range("B8").Value  = "'2025-02-23"
$ws1.Range("D8").Value  = 50

$ws1.Range("B9").Value  = "'2025-03-02"
$ws1.Range("D9").Value  = 51

$ws1.Range("B10").Value = "'2025-03-09"
$ws1.Range("D10").Value = 51

$ws1.Range("B11").Value = "'2025-03-16"
$ws1.Range("D11").Value = 51

$ws1.Range("B12").Value = "'2025-03-23"
$ws1.Range("D12").Value = 52

$ws1.Range("B13").Value = "'2025-03-30"
$ws1.Range("D13").Value = 53

$ws1.Range("B14").Value = "'2025-04-06"
$ws1.Range("D14").Value = 53

$ws1.Range("B15").Value = "'2025-04-13"
$ws1.Range("D15").Value = 53

$ws1.Range("B16").Value = "'2025-04-20"
$ws1.Range("D16").Value = 53

$ws1.Range("B17").Value = "'2025-04-27"
$ws1.Range("D17").Value = 54

# --- Sheet "Summary": recomputed stats that depend on the shifted window ---

$ws2.Range("B2").Value  = "2024-06-02 to 2025-01-05"   # Historical Range
$ws2.Range("B3").Value  = "'10"                         # Min Sales
$ws2.Range("B4").Value  = "'63"                         # Max Sales
# B5 Mean Sales unchanged
$ws2.Range("B6").Value  = "'38"                         # Median Sales
$ws2.Range("B7").Value  = "'12"                         # Std Dev Sales
$ws2.Range("B8").Value  = "1154 units"                  # Total Historical Sales
$ws2.Range("B9").Value  = "'811"                         # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "'391"                         # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "'191"                         # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "'54"                          # Max Forecast
$ws2.Range("B13").Value = "'2025-04-27"                  # Max Forecast Week
$ws2.Range("B14").Value = "'47"                          # Min Forecast
$ws2.Range("B15").Value = "'2025-01-12"                  # Min Forecast Week
